$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.392.35"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").Value = "3.778.52"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.42"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.10"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("D7").Value = "3.777.24"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.37"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  -4.00%  "
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").Value = "4.408.60"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "3.762.04"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "67.407.04"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.86"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.96"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.23"
$ws.Range("E21").Value = "  -7.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "457.15"
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.699"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000150"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.35"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.86"
$ws.Range("E26").Value = "  -3.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.12"
$ws.Range("E27").Value = "  -5.85%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -3.43%  "
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.18"
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.17"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").Value = "3.730.06"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0996"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.137"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("E39").Value = "  -7.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.992"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.71"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.297"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.88"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "147.39"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "391.43"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.82"
$ws.Range("E50").Value = "  -8.29%  "
$ws.Range("D51").Value = "2.754.11"
$ws.Range("E51").Value = "  +1.91%  "
